# Update "想去人数" (want-to-go count) / "最低票价" (lowest price) figures
# for a handful of events across the "展览", "演出" and "全部类型" sheets.
# These numbers were refreshed from the live bilibili event pages.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibition) sheet ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 506    # id=89803  505 -> 506
$wsExpo.Range("F4").Value = 30     # id=90924  27  -> 30
$wsExpo.Range("F9").Value = 111    # id=91133  91  -> 111
$wsExpo.Range("F10").Value = 2060  # id=90908  1939 -> 2060

# --- 演出 (Performance) sheet ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 42     # id=90322  41 -> 42
$wsShow.Range("G4").Value = 64     # id=90593  56 -> 64

# --- 全部类型 (All types, union of the other sheets) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 506     # id=89803  505 -> 506
$wsAll.Range("F5").Value = 30      # id=90924  27  -> 30
$wsAll.Range("F10").Value = 111    # id=91133  91  -> 111
$wsAll.Range("F11").Value = 2060   # id=90908  1939 -> 2060
$wsAll.Range("F13").Value = 42     # id=90322  41 -> 42
$wsAll.Range("G14").Value = 64     # id=90593  56 -> 64
